# Reformat the single shared-string cell (pretty-printed JSON instead of a
# Python dict literal), and collapse the sheet back down to a single used
# cell: the old numeric placeholder in A1 is removed, and the text that used
# to live in A2 now lives in A1 (with A1's previous bold/bordered/centered
# style cleared back to the workbook default).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop both previously-used cells (value + formatting) so A1 starts clean.
$ws.Range("A1:A2").Clear()

# Write the reformatted text into A1 (now the sheet's only cell).
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You are the IT administrator at a company with 150 employees. You need to migrate the company\u2019s legacy email and videoconferencing system to a cloud-based service. The company requires a product that meets the following requirements:Includes email and videoconferencing servicesProvides vendor supportMinimizes costsYou''ve decided to use Google Workspace.Which Google Workspace edition should you recommend?",
        "ques_type": 2,
        "options": [
            "Business Starter",
            "Business Plus",
            "Essentials",
            "Enterprise Standard"
        ],
        "score": "Business Starter"
    },
    {
        "title": "You are an IT administrator at a company that uses Google Workspace for email. You''ve established a group for the sales department to collect purchase requests. However, the sales manager wants you to modify the group so that purchase requests can be assigned to team members. All members must be able to see who is responsible for each request.What is the most appropriate modification to make?",
        "ques_type": 2,
        "options": [
            "Change the group type to a dynamic group.",
            "Designate the group as a security group.",
            "Change it into a collaborative inbox.",
            "Create a separate mailbox and configure delegation."
        ],
        "score": "Change it into a collaborative inbox."
    },
    {
        "title": "You''re an IT administrator for a global company with offices in London, New York, and Paris. The Paris office is experiencing network issues, causing a decrease in internet bandwidth. You decide to alleviate the network load by restricting Google Meet video quality in the Paris office.Which actions should you take to do this most effectively?",
        "ques_type": 15,
        "options": [
            "Close UDP ports 19302\u200b\u201319309 on the Paris office firewall.",
            "Move the profiles of London and New York users to a separate organizational unit.",
            "Use Meet safety settings from Google Admin Console to set the video quality.",
            "Move the profiles of Paris users to a separate organizational unit.",
            "Use Meet video settings from Google Admin Console to set the video quality."
        ],
        "score": [
            "Move the profiles of Paris users to a separate organizational unit.",
            "Use Meet video settings from Google Admin Console to set the video quality."
        ]
    },
    {
        "title": "You''re an IT manager at a healthcare company handling highly sensitive patient data. You need to enhance the security of a selected group of users'' Google Workspace accounts by enforcing hardware multi-factor authentication devices for two-step verification. You have opened Google Admin Console, navigated to Security &gt Authentication &gt 2-Step Verification, and selected the target group. Which setting should you enable?",
        "ques_type": 2,
        "options": [
            "Methods &gt Any except verification codes via text, phone call",
            "Security Codes &gt Don\u2019t allow users to generate security codes",
            "Security Codes &gt Allow security codes without remote access",
            "Methods &gt Only Security Key"
        ],
        "score": "Methods &gt Only Security Key"
    }
]'
